$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the SDET program name cells per the data change described in the commit.
# Order matters for shared-string table layout, so write in this sequence:
# A2, A3, B3, B2
$ws.Range("A2").Value = "July27-ninjalibaries-SDET321-321"
$ws.Range("A3").Value = "July27-ninjalibaries-SDET432-432"
$ws.Range("B3").Value = "June23-ninjalibaries-SDET432-432"
$ws.Range("B2").Value = "June23-ninjalibaries-SDET321-321"
